$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Status column updates (Task 1 & Task 2 now Complete) ---
$ws.Range("B2").Value = "Complete"
$ws.Range("B3").Value = "Complete"

# --- Date Last Updated column (new column D) ---
# Enter the date as literal text (not an auto-converted date serial):
# build it via a formula, then paste back as a value.
$ws.Range("D2:D3").Formula = '="01/30/2024"'
$ws.Range("D2:D3").Copy()
$ws.Range("D2:D3").PasteSpecial(-4163)

# --- New header cell ---
$ws.Range("D1").Value = "Date Last Updated"

# --- Notes for Task 2 ---
$ws.Range("C3").Value = "using max-width min-width on css to make it so"

# --- Flesh out Task 3-6 descriptions ---
$ws.Range("A4").Value = "Task 3: Login Page"
$ws.Range("A5").Value = "Task 4: Establish DB Env"
$ws.Range("A6").Value = "Task 5: DB:Schema"
$ws.Range("A7").Value = "Task 6: Create Necessary DB Tables"

$ws.Application.CutCopyMode = $false

# --- Header row bold formatting ---
$ws.Range("A1:D1").Font.Bold = $true

# --- New column width ---
$ws.Columns.Item(4).ColumnWidth = 15.8

# --- Selection / view ---
$ws.Range("C7").Select()

# --- Page setup orientation ---
$ws.PageSetup.Orientation = 1
